$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Enabling only two test cases (rows 2 and 3) in the Notification Suite;
# disable the rest by setting Runmode (column D) to "N" for rows 4-22.
for ($row = 4; $row -le 22; $row++) {
    $ws.Cells.Item($row, 4).Value = "N"
}
